# This script updates cryptocurrency price/volume figures in the active worksheet
# to match the refreshed data pulled on 2023-02-05, per the "Updated symbol list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'334.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.51%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'44.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'6.23%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.743"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.63%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08365"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.91%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'8.852"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.97%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.950"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-4.44%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D9").Value = "'0.9472"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.36%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1241"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-2.58%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1970"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.42%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'6.90%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.04417"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'12.50%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.1068"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.69%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001295"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.63%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006067"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.38%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.484"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.33%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'4.523"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.16%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value = "'8.718"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'4.27%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-0.72%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D23").Value = "'0.04411"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.18%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001247"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.77%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004362"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.00%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001264"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'5.24%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003999"
$ws.Range("D27").Style = "Normal"
$ws.Range("D39").Value = "'0.02820"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.18%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05876"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'6.42%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007942"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.07%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1427"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.30%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.009043"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'1.18%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002149"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'0.30%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01038"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-12.70%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007235"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'3.29%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.23%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003200"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'0.30%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002276"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.20%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002107"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.23%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002007"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.23%"
$ws.Range("E51").Style = "Normal"
